# Generate Report for Handback
# For both the zh-cn and de-de sheets, row 6 (the 2ac31f83-... file) gets a
# handback that is out of date: fill in the Latest Target File, Latest
# Handback DateTime and Error Detail columns, and add a hyperlink in the
# Latest Target File cell (column I) pointing at the latest commit.

$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-08-25 22:42:27" },
    @{ Name = "de-de"; HandbackTime = "2016-08-25 22:42:34" }
)

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b0346f850ed76f00783a988d7af3966f93f49d1/e2e/2ac31f83-5836-4d9b-ab59-8568bbf00d44.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d845d3e72e851b75a50ae97a92be838a2b6980d/e2e/2ac31f83-5836-4d9b-ab59-8568bbf00d44.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b0346f850ed76f00783a988d7af3966f93f49d1/e2e/2ac31f83-5836-4d9b-ab59-8568bbf00d44.md."

foreach ($entry in $sheets) {
    $ws = $wb.Worksheets.Item($entry.Name)

    # Column P ("Error Detail") needs more room for the long message.
    $ws.Columns.Item(16).ColumnWidth = 40

    # Latest Target File (I6): a hyperlink to the latest handoff markdown file.
    $ws.Range("I6").Value = "2ac31f83-5836-4d9b-ab59-8568bbf00d44.md"
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestUrl, "", "", "2ac31f83-5836-4d9b-ab59-8568bbf00d44.md")

    # Latest Handback File (J6): same xlf file name as the handoff file (G6).
    $ws.Range("J6").Value = $ws.Range("G6").Value2

    # Latest Handback DateTime (K6).
    $ws.Range("K6").Value = $entry.HandbackTime

    # Error Detail (P6): version mismatch message.
    $ws.Range("P6").Value = $errorDetail
}
